{"js": "// Apply the \"music and art assets\" edit:\n// 1. Append a sentence about bonus money for shooting asteroids to the\n//    \"2nd weapon system ... missiles not cheap.\" paragraph, and move the\n//    _GoBack bookmark there (to the end of that paragraph).\n// 2. Remove the two TODO paragraphs \"Make crashed ship obstacle sprite\"\n//    and \"Machine gun/weaponry sprites\".\n// 3. Add a new TODO paragraph \"Make game look nice when maximized\" right\n//    after the \"All audio (...)\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraphs we need to act on by their (stable) text content.\nlet missileParaIndex = -1;\nlet crashedShipParaIndex = -1;\nlet machineGunParaIndex = -1;\nlet allAudioParaIndex = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (missileParaIndex === -1 && t.indexOf(\"Moderate expense, missiles not cheap.\") !== -1) {\n    missileParaIndex = i;\n  }\n  if (crashedShipParaIndex === -1 && t.indexOf(\"Make crashed ship obstacle sprite\") !== -1) {\n    crashedShipParaIndex = i;\n  }\n  if (machineGunParaIndex === -1 && t.indexOf(\"Machine gun/weaponry sprites\") !== -1) {\n    machineGunParaIndex = i;\n  }\n  if (allAudioParaIndex === -1 && t.indexOf(\"All audio (collision sound effects\") !== -1) {\n    allAudioParaIndex = i;\n  }\n}\n\nif (missileParaIndex === -1) throw new Error(\"Could not find the 'missiles not cheap' paragraph\");\nif (crashedShipParaIndex === -1) throw new Error(\"Could not find the 'Make crashed ship obstacle sprite' paragraph\");\nif (machineGunParaIndex === -1) throw new Error(\"Could not find the 'Machine gun/weaponry sprites' paragraph\");\nif (allAudioParaIndex === -1) throw new Error(\"Could not find the 'All audio' paragraph\");\n\n// --- Change 2: remove the two sprite TODO paragraphs ---\nparagraphs.items[crashedShipParaIndex].delete();\nparagraphs.items[machineGunParaIndex].delete();\n\n// --- Change 3: add the new TODO paragraph after \"All audio (...)\" ---\nparagraphs.items[allAudioParaIndex].insertParagraph(\n  \"Make game look nice when maximized\",\n  Word.InsertLocation.after\n);\n\n// --- Change 1: append the bonus-money sentence + move the _GoBack bookmark ---\n// The _GoBack bookmark currently lives in the \"All audio\" paragraph; remove\n// it there (document-wide, name is unique) before re-inserting it at its new\n// location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst missilePara = paragraphs.items[missileParaIndex];\nmissilePara.insertText(\n  \" Get bonus money for shooting down asteroids with missiles (Earth Defense!)\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\nconst endOfMissilePara = missilePara.getRange(Word.RangeLocation.end);\nendOfMissilePara.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Apply the \"music and art assets\" edit:\n# 1. Append a sentence about bonus money for shooting asteroids to the\n#    \"2nd weapon system ... missiles not cheap.\" paragraph, and move the\n#    _GoBack bookmark there (to the very end of that paragraph).\n# 2. Remove the two TODO paragraphs \"Make crashed ship obstacle sprite\"\n#    and \"Machine gun/weaponry sprites\".\n# 3. Add a new TODO paragraph \"Make game look nice when maximized\" right\n#    after the \"All audio (...)\" paragraph.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphByText($searchText) {\n    $r = $d.Content\n    $find = $r.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Could not find text: $searchText\"\n    }\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Start -le $r.Start -and $p.Range.End -ge $r.End) {\n            return $p\n        }\n    }\n    throw \"Could not locate containing paragraph for: $searchText\"\n}\n\n# --- Change 2: remove the two sprite TODO paragraphs ---\n$pCrashedShip = Get-ParagraphByText(\"Make crashed ship obstacle sprite\")\n$pCrashedShip.Range.Delete()\n\n$pMachineGun = Get-ParagraphByText(\"Machine gun/weaponry sprites\")\n$pMachineGun.Range.Delete()\n\n# --- Change 3: add the new TODO paragraph after \"All audio (...)\" ---\n$pAllAudio = Get-ParagraphByText(\"All audio (collision sound effects\")\n$allAudioIndex = $pAllAudio.Index\n$endOfAllAudio = $pAllAudio.Range\n[void]$endOfAllAudio.MoveEnd(1, -1)\n[void]$endOfAllAudio.Collapse(0)\n$endOfAllAudio.InsertParagraphAfter()\n\n# The freshly inserted (empty) paragraph is the very next one; fill it in.\n$newPara = $d.Paragraphs.Item($allAudioIndex + 1)\n$newRange = $newPara.Range\n[void]$newRange.MoveEnd(1, -1)\n$newRange.Text = \"Make game look nice when maximized\"\n\n# --- Change 1: append the bonus-money sentence + move the _GoBack bookmark ---\n# The _GoBack bookmark currently lives in the \"All audio\" paragraph; remove it\n# (its name is unique document-wide) before re-inserting it at its new spot.\n$goBack = $d.Bookmarks.Item(\"_GoBack\")\n$goBack.Delete()\n\n$pMissile = Get-ParagraphByText(\"Moderate expense, missiles not cheap.\")\n$missileRange = $pMissile.Range\n[void]$missileRange.MoveEnd(1, -1)\n$missileRange.InsertAfter(\" Get bonus money for shooting down asteroids with missiles (Earth Defense!)\")\n\n# Re-fetch the paragraph range (InsertAfter grew it) and find its new end,\n# excluding the trailing paragraph mark.\n$pMissile2 = Get-ParagraphByText(\"Moderate expense, missiles not cheap.\")\n$missileRange2 = $pMissile2.Range\n[void]$missileRange2.MoveEnd(1, -1)\n$endPos = $missileRange2.End\n\n# Work around the engine only placing Bookmarks.Add correctly for\n# non-collapsed ranges: insert a 1-char placeholder at the target spot,\n# bookmark around it, then delete the placeholder. The bookmark collapses\n# to a proper zero-width _GoBack marker at exactly that position.\n$placeholderRange = $d.Range($endPos - 1, $endPos)\n$placeholderRange.InsertAfter(\"X\")\n$markerRange = $d.Range($endPos, $endPos + 1)\n$d.Bookmarks.Add(\"_GoBack\", $markerRange) | Out-Null\n$deleteRange = $d.Range($endPos, $endPos + 1)\n$deleteRange.Delete()\n"}
